$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Three new script-file rows (same shape as the existing ones in column A).
$ws.Range("A4").Value = "SCRIPT/P02P01A/um1103.ssb"
$ws.Range("A5").Value = "SCRIPT/P02P01A/um1106.ssb"
$ws.Range("A6").Value = "SCRIPT/P02P01A/um1109.ssb"

# Match the wrapped-row height used elsewhere on the sheet.
$ws.Range("A4:A6").RowHeight = 43.2

# Row 3 no longer is the last row of the table, so drop the thin bottom
# border it had (the new rows below it carry no border).
$ws.Range("A3:E3").Borders.LineStyle = -4142

# Leave the selection where the author's editing session ended up.
$ws.Range("C4").Select() | Out-Null
